# Generate Report for Handoff
#
# Swaps the two tracked files' report rows:
#   - "a40c3310-6a44-4511-8405-73d0d835b0fe.md" moves into the row that used
#     to hold "8bdff04d-2ac5-4fa2-8756-f40df755edd9.md" (and vice versa) on
#     every sheet, since the report tool re-sorted the rows.
#   - The 8bdff04d file's status flips from "In Translation" to
#     "Ready for handoff" (zh-cn and the Overview roll-up), its zh-cn
#     handoff timestamp advances, and its Priority becomes "mt".
#   - The a40c3310 file's zh-cn status also flips to "Ready for handoff".

$wb = $excel.ActiveWorkbook

# Helper: write a literal "True"/"False" text label into a cell without the
# engine coercing it into a real boolean cell type (t="b"). A leading
# apostrophe forces text, and resetting the style back to Normal drops the
# "quote prefix" flag that Excel would otherwise stamp on the cell.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "a40c3310-6a44-4511-8405-73d0d835b0fe.md"
$ov.Range("B2").Value = "e2e\a40c3310-6a44-4511-8405-73d0d835b0fe.md"
$ov.Range("C2").Value = ".md"
$ov.Range("E2").Value = "In Translation"
$ov.Range("F2").Value = "In Translation"
$ov.Range("G2").Value = "2016-08-13 20:17:42"

$ov.Range("A3").Value = "8bdff04d-2ac5-4fa2-8756-f40df755edd9.md"
$ov.Range("B3").Value = "e2e\8bdff04d-2ac5-4fa2-8756-f40df755edd9.md"
$ov.Range("C3").Value = ".md"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-13 20:18:34"

foreach ($h in $ov.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$2') {
        $h.TextToDisplay = "e2e\a40c3310-6a44-4511-8405-73d0d835b0fe.md"
    }
    if ($addr -eq '$B$3') {
        $h.TextToDisplay = "e2e\8bdff04d-2ac5-4fa2-8756-f40df755edd9.md"
    }
}

$ov.Columns.Item(5).EntireColumn.AutoFit() | Out-Null
$ov.Columns.Item(6).EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "a40c3310-6a44-4511-8405-73d0d835b0fe.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("D2").Value = "e2e"
$zh.Range("E2").Value = "ht"
Set-TextValue $zh.Range("F2") "False"
$zh.Range("G2").Value = "a40c3310-6a44-4511-8405-73d0d835b0fe.aecf525dc11b0f093f604912985ceda5e3771253.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-13 20:17:34"
$zh.Range("K2").Value = "0001-01-01 00:00:00"
Set-TextValue $zh.Range("M2") "True"
Set-TextValue $zh.Range("O2") "False"

$zh.Range("A3").Value = "8bdff04d-2ac5-4fa2-8756-f40df755edd9.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "e2e"
$zh.Range("E3").Value = "mt"
Set-TextValue $zh.Range("F3") "False"
$zh.Range("G3").Value = "8bdff04d-2ac5-4fa2-8756-f40df755edd9.723173f88eafcd10730c03f4494a2eee648f3cb2.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-13 20:18:26"
$zh.Range("K3").Value = "0001-01-01 00:00:00"
Set-TextValue $zh.Range("M3") "True"
Set-TextValue $zh.Range("O3") "False"

foreach ($h in $zh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "a40c3310-6a44-4511-8405-73d0d835b0fe.md"
    }
    if ($addr -eq '$A$3') {
        $h.TextToDisplay = "8bdff04d-2ac5-4fa2-8756-f40df755edd9.md"
    }
}

$zh.Columns.Item(3).EntireColumn.AutoFit() | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "a40c3310-6a44-4511-8405-73d0d835b0fe.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "In Translation"
$de.Range("D2").Value = "e2e"
$de.Range("E2").Value = "ht"
Set-TextValue $de.Range("F2") "False"
$de.Range("G2").Value = "a40c3310-6a44-4511-8405-73d0d835b0fe.aecf525dc11b0f093f604912985ceda5e3771253.de-de.xlf"
$de.Range("H2").Value = "2016-08-13 20:17:42"
$de.Range("K2").Value = "0001-01-01 00:00:00"
Set-TextValue $de.Range("M2") "True"
Set-TextValue $de.Range("O2") "False"

$de.Range("A3").Value = "8bdff04d-2ac5-4fa2-8756-f40df755edd9.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "e2e"
$de.Range("E3").Value = "mt"
Set-TextValue $de.Range("F3") "False"
$de.Range("G3").Value = "8bdff04d-2ac5-4fa2-8756-f40df755edd9.723173f88eafcd10730c03f4494a2eee648f3cb2.de-de.xlf"
$de.Range("H3").Value = "2016-08-13 20:18:34"
$de.Range("K3").Value = "0001-01-01 00:00:00"
Set-TextValue $de.Range("M3") "True"
Set-TextValue $de.Range("O3") "False"

foreach ($h in $de.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = "a40c3310-6a44-4511-8405-73d0d835b0fe.md"
    }
    if ($addr -eq '$A$3') {
        $h.TextToDisplay = "8bdff04d-2ac5-4fa2-8756-f40df755edd9.md"
    }
}

$de.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
